# Update Name of Algo
# Apply updated RandomForest imputation results to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B3").Value = 6.173700000000005
$ws.Range("C3").Value = -11.8505
$ws.Range("A12").Value = -21.65990000000001
$ws.Range("B14").Value = 6.404099999999997
$ws.Range("B26").Value = 3.850700000000002
$ws.Range("C30").Value = -12.94709999999999
$ws.Range("B31").Value = 4.491200000000003
$ws.Range("A32").Value = -21.2776
$ws.Range("B35").Value = 9.187200000000006
$ws.Range("A36").Value = -19.5135
$ws.Range("B37").Value = 8.989900000000004
$ws.Range("A38").Value = -19.5382
$ws.Range("C44").Value = -13.59899999999999
$ws.Range("B45").Value = 5.598000000000003
$ws.Range("A46").Value = -21.74020000000001
$ws.Range("A54").Value = -21.85859999999998
$ws.Range("A55").Value = -22.39670000000001
$ws.Range("B57").Value = 5.114099999999997
$ws.Range("C58").Value = -13.417
$ws.Range("A67").Value = -21.43729999999997
$ws.Range("A69").Value = -21.58569999999999
$ws.Range("A72").Value = -21.8537
$ws.Range("C84").Value = -13.80799999999999
$ws.Range("C89").Value = -11.251
$ws.Range("A91").Value = -21.50610000000001
$ws.Range("C91").Value = -11.0225
$ws.Range("C92").Value = -11.3673
$ws.Range("A99").Value = -20.35759999999999
$ws.Range("B100").Value = 5.148199999999998
$ws.Range("B102").Value = 8.113500000000004
$ws.Range("C102").Value = -13.0024
